# Apply the "Fruta / hortaliza, semanal" update:
# Insert two new weekly price records at rows 46-47, shifting the
# previously existing rows 46-68 down to rows 48-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 46 (this pushes old
# rows 46-68 down to 48-70, matching the diff's row-shift pattern).
$ws.Range("A46:R47").Insert()

# --- New row 46 ---
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44680
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = 100112021
$ws.Cells.Item(46, 7).Value = "Ají"
$ws.Cells.Item(46, 8).Value = "Inferno"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 120
$ws.Cells.Item(46, 11).Value = 25000
$ws.Cells.Item(46, 12).Value = 26000
$ws.Cells.Item(46, 13).Value = 25500
$ws.Cells.Item(46, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 1700
$ws.Cells.Item(46, 17).Value = 15
$ws.Cells.Item(46, 18).Value = "Hortaliza"

# --- New row 47 ---
$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(47, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value = 44680
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(47, 6).Value = 100112021
$ws.Cells.Item(47, 7).Value = "Ají"
$ws.Cells.Item(47, 8).Value = "Inferno"
$ws.Cells.Item(47, 9).Value = "Segunda"
$ws.Cells.Item(47, 10).Value = 130
$ws.Cells.Item(47, 11).Value = 21000
$ws.Cells.Item(47, 12).Value = 22000
$ws.Cells.Item(47, 13).Value = 21500
$ws.Cells.Item(47, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 1433
$ws.Cells.Item(47, 17).Value = 15
$ws.Cells.Item(47, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
